$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(15, 1).Value = 9858.19
$ws.Cells.Item(15, 2).Value = 9912.7099999999991
$ws.Cells.Item(15, 3).Value = 78.48
$ws.Cells.Item(15, 4).Value = 78.05
$ws.Cells.Item(15, 5).Value = $false
$ws.Cells.Item(15, 6).Value = -0.55000000000000004
$ws.Cells.Item(15, 7).Value = 42624.611145833333
$ws.Cells.Item(15, 8).Value = $false
